$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 11: add the missing "total qty" cells next to each "total amt" cell ---
# The existing amt cells (F11, H11, J11, L11) keep their column/value; we only
# need to populate the previously-empty qty cells that precede them
# (E11, G11, I11, K11) and restyle the whole row to right-align.
#
#   E11=total_open_qty        F11=total_open_amt        (unchanged)
#   G11=total_period_in_qty   H11=total_period_in_amt   (unchanged)
#   I11=total_period_out_qty  J11=total_period_out_amt  (unchanged)
#   K11=total_close_qty       L11=total_close_amt       (unchanged)

# F11/H11 carry the "no border" bold style family -- reuse it for the new
# E/G/I "qty" cells, then right-align the whole block (E,F,G,H,I).
$ws.Range("F11").Copy($ws.Range("E11"))
$ws.Range("F11").Copy($ws.Range("G11"))
$ws.Range("F11").Copy($ws.Range("I11"))

$ws.Range("E11").Value2 = "{{currency total_open_qty}}"
$ws.Range("G11").Value2 = "{{currency total_period_in_qty}}"
$ws.Range("I11").Value2 = "{{currency total_period_out_qty}}"

$ws.Range("E11:I11").HorizontalAlignment = -4152

# J11/L11 carry the "with border" bold style family -- reuse it for the new
# K11 "qty" cell, then right-align the whole block (J,K,L) while keeping
# their existing vertical-centered alignment.
$ws.Range("J11").Copy($ws.Range("K11"))

$ws.Range("K11").Value2 = "{{currency total_close_qty}}"

$ws.Range("J11:L11").HorizontalAlignment = -4152

# --- Sheet view: selection moved to G12 ---
$ws.Range("G12").Select() | Out-Null
